# Auto-generated edit script: update cryptos list values per commit
# "Updated cryptos list on Thu Aug 17 12:46:45 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.632.36"
$ws.Range("E2").Value = "  -2.20%  "
# Row 3
$ws.Range("D3").Value = "1.796.49"
$ws.Range("E3").Value = "  -1.96%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.44"
$ws.Range("E5").Value = "  -1.17%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5915"
$ws.Range("E6").Value = "  -1.45%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.09%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2787"
$ws.Range("E8").Value = "  +0.69%  "
# Row 9
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.37"
$ws.Range("E9").Value = "  -0.25%  "
# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06777"
$ws.Range("E10").Value = "  -3.05%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07554"
$ws.Range("E11").Value = "  -1.22%  "
# Row 12
$ws.Range("D12").Value = "1.795.57"
$ws.Range("E12").Value = "  -1.80%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.805"
$ws.Range("E13").Value = "  +0.45%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6169"
$ws.Range("E14").Value = "  -1.38%  "
# Row 15
$ws.Range("D15").Value = "2.039.49"
$ws.Range("E15").Value = "  -1.79%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "75.59"
$ws.Range("E16").Value = "  -4.10%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008943"
$ws.Range("E17").Value = "  -7.68%  "
# Row 18
$ws.Range("D18").Value = "28.608.04"
$ws.Range("E18").Value = "  -2.06%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.420"
$ws.Range("E19").Value = "  -6.12%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  -0.05%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "209.99"
$ws.Range("E21").Value = "  -6.12%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.49"
$ws.Range("E22").Value = "  -1.03%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.854"
$ws.Range("E23").Value = "  -0.61%  "
# Row 24
$ws.Range("E24").Value = "  -0.02%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.79"
$ws.Range("E25").Value = "  -2.53%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.194"
$ws.Range("E26").Value = "  +2.76%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1265"
$ws.Range("E27").Value = "  -2.70%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.48"
$ws.Range("E28").Value = "  -0.64%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.414"
$ws.Range("E29").Value = "  -2.69%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06283"
$ws.Range("E30").Value = "  -6.70%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.424"
$ws.Range("E31").Value = "  -1.51%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.800"
$ws.Range("E32").Value = "  -1.00%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.791"
$ws.Range("E33").Value = "  +0.54%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.742"
$ws.Range("E34").Value = "  +0.98%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.051"
$ws.Range("E35").Value = "  -4.37%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6433"
$ws.Range("E36").Value = "  -0.13%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.504"
$ws.Range("E37").Value = "  -1.64%  "
# Row 38
$ws.Range("E38").Value = "  -1.02%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01702"
$ws.Range("E39").Value = "  -2.43%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.360"
$ws.Range("E40").Value = "  -2.07%  "
# Row 41
$ws.Range("D41").Value = "1.150.01"
$ws.Range("E41").Value = "  -4.39%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8764"
$ws.Range("E42").Value = "  -2.90%  "
# Row 43
$ws.Range("E43").Value = "  -0.14%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.32"
$ws.Range("E44").Value = "  -0.03%  "
# Row 45
$ws.Range("D45").Value = "1.945.94"
$ws.Range("E45").Value = "  -1.83%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.22"
$ws.Range("E46").Value = "  -3.43%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000109"
$ws.Range("E47").Value = "  -3.47%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.593"
$ws.Range("E48").Value = "  +0.85%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.348"
$ws.Range("E49").Value = "  -1.68%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05465"
$ws.Range("E50").Value = "  -0.81%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4478"
$ws.Range("E51").Value = "  -1.75%  "
